$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new rows before the old "empty row / total" block (old rows 24-26) ---
# This shifts old row24->29, old row25->30, old row26 (Total)->31
$ws.Rows("24:28").Insert()

# Copy formatting (styles/borders) from row 23 down into the newly inserted rows
# so the blank rows keep the same look as the rest of the data rows.
$ws.Range("A23:F23").Copy()
$ws.Range("A24:F30").PasteSpecial(-4122)   # xlPasteFormats

# The insert leaves a left-over "=C-B" formula (evaluating to 0) in every row of
# the shifted block (now rows 24-29); the target only wants that formula to
# remain on the very last inserted row (row 30) - the rest (24-29) go back to
# being truly blank placeholder rows.
$ws.Range("D24:D29").ClearContents()

# --- Fill in the new diary entry (row 24): 14.11.18, 10:15 -> 20:15 ---
$ws.Range("A24").Value = "14.11.18"
$ws.Range("B24").Value = 0.42708333333333331
$ws.Range("C24").Value = 0.84375
$ws.Range("D24").Formula = "=C24-B24"
$ws.Range("E24").Value = "-Output Options"

# --- Update the remark text for the existing rows 22 and 23 ---
# (Order matters for shared-string table allocation: TriggerEvent is added
# before Preparing OutputOptions so new indices land the same as the target.)
$ws.Range("E23").Value = "-Implementation of TriggerEvent"
$ws.Range("E22").Value = "-Preparing OutputOptions"

# --- Update view state: scrolled down a bit, selection sitting on E36 ---
$ws.Range("E36").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 16
